$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Littleneck clam (Leukoma staminea)"
$ws.Range("A9").Value = "Cockle (Clinocardium nuttallii)"

$ws.Range("A10").Select()
